# Tesla Valuation v1.0.xlsx - Continuing to update data for Relative Valuation
# Adds an "EV/EBITDA 2021E" (D) data column and an "EV/EBITDA 2021E x1.2" (F)
# derived column to the peer-comparison table (rows 58-67), wires up the
# corresponding multiple columns (H/I/J), and updates the Average/Median
# summary rows (71/72) to cover the new D-based multiple.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58-67: add column D (new EBITDA estimate), matching column C's format ---
$ws.Range("C58:C67").Copy()
$ws.Range("D58:D67").PasteSpecial(-4122)

# --- Prime H58:J67 with the same format currently used by G58 (s=4 "Accent1 x" style) ---
$ws.Range("G58").Copy()
$ws.Range("H58:J58").PasteSpecial(-4122)
$ws.Range("H59:J67").PasteSpecial(-4122)

# --- New column D values ($m) ---
$ws.Range("D58").Value = 333420000000
$ws.Range("D59").Value = 439000000000
$ws.Range("D60").Value = 226100000000
$ws.Range("D61").Value = 107550000000
$ws.Range("D62").Value = 30050000000
$ws.Range("D63").Value = 164190000000
$ws.Range("D64").Value = 13790000000
$ws.Range("D65").Value = 108560000000
$ws.Range("D66").Value = 15510000000
$ws.Range("D67").Value = 72450000000

# --- New column F58: derived EBITDA figure (E58 * 1.2), default "Normal" style ---
$ws.Range("F58").Formula = "=E58*1.2"
$ws.Range("F58").Style = "Normal"

# --- Row 58 multiple formulas: existing EV/E multiple moves from H to I,
#     H58 becomes the new EV/D multiple, J58 is the new EV/F multiple ---
$ws.Range("I58").Formula = "=ROUND((B58/E58),2)&""x"""
$ws.Range("H58").Formula = "=ROUND((B58/D58),2)&""x"""
$ws.Range("J58").Formula = "=ROUND((B58/F58),2)&""x"""

# --- Rows 59-67: fill I/J (new shared formulas) before H so the shared-formula
#     index allocation lines up with a from-scratch edit (H reuses its existing
#     shared group; I/J are brand new groups) ---
$ws.Range("I59:I67").Formula = "=ROUND((B59/E59),2)&""x"""
$ws.Range("J59:J67").Formula = "=ROUND((B59/F59),2)&""x"""
$ws.Range("H59:H67").Formula = "=ROUND((B59/D59),2)&""x"""

# --- Row 71 (Average) / Row 72 (Median): extend to the new D-based multiple,
#     matching the "Output, right-aligned" style (s=39) used elsewhere (e.g. H43) ---
$ws.Range("H43").Copy()
$ws.Range("G71:J71").PasteSpecial(-4122)
$ws.Range("G72:J72").PasteSpecial(-4122)

$ws.Range("G71").Formula = "=ROUND(AVERAGE((B58/C58),(B59/C59),(B60/C60),(B61/C61),(B62/C62),(B63/C63),(B64/C64),(B65/C65),(B66/C66),(B67/C67)),2)&""x"""
$ws.Range("H71").Formula = "=ROUND(AVERAGE((B58/D58),(B59/D59),(B60/D60),(B61/D61),(B62/D62),(B63/D63),(B64/D64),(B65/D65),(B66/D66),(B67/D67)),2)&""x"""

$ws.Range("G72").Formula = "=ROUND(MEDIAN((B58/C58),(B59/C59),(B60/C60),(B61/C61),(B62/C62),(B63/C63),(B64/C64),(B65/C65),(B66/C66),(B67/C67)),2)&""x"""
$ws.Range("H72").Formula = "=ROUND(MEDIAN((B58/D58),(B59/D59),(B60/D60),(B61/D61),(B62/D62),(B63/D63),(B64/D64),(B65/D65),(B66/D66),(B67/D67)),2)&""x"""

# --- Update the sheet's active selection (view state) ---
$ws.Range("F81").Select() | Out-Null
